$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.226335644721985
$ws.Range("B1").Value = 2.704195499420166
$ws.Range("C1").Value = 4.289385795593262
$ws.Range("D1").Value = 2.117380857467651
$ws.Range("E1").Value = 1.155708789825439
